$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 (client record #15) was previously blank in columns B:G.
# Fill in the new client record that was uploaded: VANDER LUIS.
$ws.Range("B16").Value = "VANDER LUIS"
$ws.Range("C16").Value = "46faf0b02e80945f4d911de265fda99a"
$ws.Range("D16").Value = 44851
$ws.Range("E16").Value = 365
$ws.Range("F16").Value = "-"
$ws.Range("G16").Value = "VENDA 12 (17/10)"

# The EMAIL / CONFIRMADO columns (F:G) on this row still carried the
# "unused template row" formatting; restore the formatting used by the
# other populated data rows (copy format from the row above).
$ws.Range("F15:G15").Copy()
$ws.Range("F16:G16").PasteSpecial(-4122)
